$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update energy consumption values (priority-queue based computation results)
$ws.Range("B2").Value = 2.639363959198433
$ws.Range("C2").Value = 0.5245111508144374
$ws.Range("B3").Value = 2.853365055017402
$ws.Range("C3").Value = 1.405007082596767
$ws.Range("B4").Value = 4.647594031840149
$ws.Range("C4").Value = 1.808882914226028
$ws.Range("B5").Value = 5.171465639361396
$ws.Range("C5").Value = 2.844151432931404
$ws.Range("B6").Value = 8.979334367206684
$ws.Range("C6").Value = 3.271682412749235
$ws.Range("B7").Value = 9.083470661467837
$ws.Range("C7").Value = 3.954090033301189
$ws.Range("B8").Value = 9.944585512194319
$ws.Range("C8").Value = 4.545064267924576
$ws.Range("B9").Value = 10.12121967950107
$ws.Range("C9").Value = 5.088475944106722
$ws.Range("B10").Value = 15.90434358016013
$ws.Range("C10").Value = 5.605124055639929
$ws.Range("B11").Value = 17.39079771196148
$ws.Range("C11").Value = 6.088141083183289
$ws.Range("B12").Value = 20.11324556069098
$ws.Range("C12").Value = 6.530572277620994
$ws.Range("B13").Value = 20.22147741143464
$ws.Range("C13").Value = 6.987227478839909
$ws.Range("B14").Value = 20.87287821065503
$ws.Range("C14").Value = 7.619186943215801
$ws.Range("B15").Value = 25.61636000599569
$ws.Range("C15").Value = 8.148935862263054
$ws.Range("B16").Value = 25.7269299242515
$ws.Range("C16").Value = 8.848526622804803
$ws.Range("B17").Value = 31.25395584817196
$ws.Range("C17").Value = 9.615584012597576
$ws.Range("B18").Value = 31.74284695724332
$ws.Range("C18").Value = 10.26177825220126
$ws.Range("B19").Value = 36.08229755404678
$ws.Range("C19").Value = 10.82589530986362
$ws.Range("B20").Value = 41.41016176168067
$ws.Range("C20").Value = 11.44125378171042
$ws.Range("B21").Value = 41.48374633029086
$ws.Range("C21").Value = 12.00005557618583
$ws.Range("B22").Value = 42.63742498485434
$ws.Range("C22").Value = 13.13579430632588
$ws.Range("B23").Value = 42.74402238975347
$ws.Range("C23").Value = 13.64421659677107
$ws.Range("B24").Value = 43.42090911939952
$ws.Range("C24").Value = 14.32685174381981
$ws.Range("B25").Value = 45.07788188771953
$ws.Range("C25").Value = 14.84887929002187
$ws.Range("B26").Value = 45.16067945582189
$ws.Range("C26").Value = 15.45692313809053
$ws.Range("B27").Value = 53.08792924880444
$ws.Range("C27").Value = 16.06711457664255
$ws.Range("B28").Value = 53.24762897162098
$ws.Range("C28").Value = 16.68520652451175
$ws.Range("B29").Value = 59.41063985236721
$ws.Range("C29").Value = 17.11651197892609
$ws.Range("B30").Value = 59.94974655683151
$ws.Range("C30").Value = 17.80193596325687
$ws.Range("B31").Value = 60.04475760209485
$ws.Range("C31").Value = 18.26689754820427
$ws.Range("B32").Value = 65.10136962377759
$ws.Range("C32").Value = 18.88357133094462
$ws.Range("B33").Value = 68.38850227086169
$ws.Range("C33").Value = 19.48112216972896
$ws.Range("B34").Value = 69.95416877722057
$ws.Range("C34").Value = 19.89180260296246
$ws.Range("B35").Value = 70.02985250269613
$ws.Range("C35").Value = 20.48642420688477
$ws.Range("B36").Value = 77.94807947417294
$ws.Range("C36").Value = 21.35267482140528
$ws.Range("B37").Value = 78.44755802865359
$ws.Range("C37").Value = 21.89578653704104
$ws.Range("B38").Value = 79.63154079910001
$ws.Range("C38").Value = 22.33627492104021
$ws.Range("B39").Value = 79.78029813082594
$ws.Range("C39").Value = 22.98100002236703
$ws.Range("B40").Value = 83.62277638944475
$ws.Range("C40").Value = 23.57428073564518
$ws.Range("B41").Value = 84.03886701799907
$ws.Range("C41").Value = 24.10004514169572
$ws.Range("B42").Value = 86.95867171411277
$ws.Range("C42").Value = 24.77433304287731
$ws.Range("B43").Value = 87.07739986776369
$ws.Range("C43").Value = 25.39711635319186
$ws.Range("B44").Value = 87.25917559543797
$ws.Range("C44").Value = 26.09585931944289
$ws.Range("B45").Value = 87.39435185479994
$ws.Range("C45").Value = 26.7749029301241
$ws.Range("B46").Value = 96.2571192206842
$ws.Range("C46").Value = 27.18475680950223

# Remove now-unused trailing rows 47-50 (data series shortened)
$ws.Range("A47:C50").EntireRow.Delete()

